$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.171.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.905.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5235"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3770"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07250"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9041"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08544"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.913.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "96.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.295"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008644"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.195.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.079"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.149.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.439"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.313"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.748"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.922"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.817"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09306"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8058"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05058"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("E35").Value = "  +0.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.442"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.64%  "
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.617"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5733"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02002"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.145"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.645"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "116.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1518"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4868"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.621"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.22%  "
